# UK removed in RP3
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ERT_SU_CZ")
$ws2 = $wb.Worksheets.Item("Change Log")

# 1. Update header label
$ws1.Range("A6").Value = "SES Area (RP3)"

# 2. Delete the "United Kingdom" row (row 36) entirely - this shifts formulas/rows up
$ws1.Rows.Item(36).Delete()

# 3. Fill in the Change Log entry
$ws2.Range("A2").Value = 44351
$ws2.Range("B2").Value = "UK"
$ws2.Range("C2").Value = 2020
$ws2.Range("D2").Value = "UK removed from RP3 area"
